$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$fmtSrc = $ws.Range("C2")

$ws.Range("E2").Value = '2026-02-19 23:18:22'
$ws.Range("I2").Value = '4.6 mm'
$ws.Range("E3").Value = '2026-02-19 23:18:25'
$ws.Range("I3").Value = '6.7 mm'
$ws.Range("O3").Value = '-6.2 °C'
$ws.Range("E4").Value = '2026-02-19 23:18:27'
$ws.Range("J4").Value = '1010.5 hPa'
$ws.Range("O4").Value = '11.7 °C'
$ws.Range("E5").Value = '2026-02-19 23:18:30'
$ws.Range("I5").Value = '8.4 mm'
$ws.Range("E6").Value = '2026-02-19 23:18:32'
$ws.Range("J6").Value = '1010.6 hPa'
$ws.Range("E7").Value = '2026-02-19 23:18:35'
$ws.Range("J7").Value = '1011.6 hPa'
$ws.Range("K7").Value = '13.6 MJ/m2'
$ws.Range("O7").Value = '13.7 °C'
$ws.Range("E8").Value = '2026-02-19 23:18:37'
$ws.Range("J8").Value = '1011.4 hPa'
$ws.Range("O8").Value = '9.7 °C'
$ws.Range("E9").Value = '2026-02-19 23:18:40'
$ws.Range("L9").Value = '41.0 km/h - 239º 22:50 TU'
$ws.Range("E10").Value = '2026-02-19 23:18:42'
$ws.Range("N10").Value = '2.9 °C 22:53 TU'
$ws.Range("O10").Value = '9.7 °C'
$ws.Range("E11").Value = '2026-02-19 23:18:45'
$c = $ws.Range("H11")
$c.NumberFormat = "@"
$c.Value = '60%'
$fmtSrc.Copy()
$c.PasteSpecial(-4122)

$ws.Range("O11").Value = '6.1 °C'
$ws.Range("E12").Value = '2026-02-19 23:18:47'
$ws.Range("E13").Value = '2026-02-19 23:18:50'
$c = $ws.Range("H13")
$c.NumberFormat = "@"
$c.Value = '60%'
$fmtSrc.Copy()
$c.PasteSpecial(-4122)

$ws.Range("J13").Value = '1012.0 hPa'
$ws.Range("E14").Value = '2026-02-19 23:18:52'
$ws.Range("O14").Value = '13.0 °C'
$ws.Range("E15").Value = '2026-02-19 23:18:54'
$c = $ws.Range("H15")
$c.NumberFormat = "@"
$c.Value = '74%'
$fmtSrc.Copy()
$c.PasteSpecial(-4122)

$ws.Range("O15").Value = '9.9 °C'
$ws.Range("E16").Value = '2026-02-19 23:18:57'
$c = $ws.Range("H16")
$c.NumberFormat = "@"
$c.Value = '76%'
$fmtSrc.Copy()
$c.PasteSpecial(-4122)

$ws.Range("I16").Value = '12.1 mm'
$ws.Range("E17").Value = '2026-02-19 23:18:59'
$c = $ws.Range("H17")
$c.NumberFormat = "@"
$c.Value = '78%'
$fmtSrc.Copy()
$c.PasteSpecial(-4122)

$ws.Range("O17").Value = '0.4 °C'
$ws.Range("E18").Value = '2026-02-19 23:19:02'
$c = $ws.Range("H18")
$c.NumberFormat = "@"
$c.Value = '62%'
$fmtSrc.Copy()
$c.PasteSpecial(-4122)

$ws.Range("J18").Value = '1010.8 hPa'
$ws.Range("N18").Value = '3.9 °C 22:53 TU'
$ws.Range("O18").Value = '11.1 °C'
$ws.Range("E19").Value = '2026-02-19 23:19:04'
$ws.Range("E20").Value = '2026-02-19 23:19:07'
$c = $ws.Range("H20")
$c.NumberFormat = "@"
$c.Value = '87%'
$fmtSrc.Copy()
$c.PasteSpecial(-4122)

$ws.Range("E21").Value = '2026-02-19 23:19:09'
$ws.Range("J21").Value = '1012.0 hPa'
$ws.Range("L21").Value = '64.1 km/h - 25º 22:53 TU'
$ws.Range("E22").Value = '2026-02-19 23:19:11'
$c = $ws.Range("H22")
$c.NumberFormat = "@"
$c.Value = '81%'
$fmtSrc.Copy()
$c.PasteSpecial(-4122)

$ws.Range("I22").Value = '1.6 mm'
$ws.Range("L22").Value = '127.8 km/h - 353º 22:39 TU'
$ws.Range("E23").Value = '2026-02-19 23:19:14'
$ws.Range("I23").Value = '12.5 mm'
$ws.Range("E24").Value = '2026-02-19 23:19:16'
$ws.Range("J24").Value = '1015.6 hPa'
$ws.Range("O24").Value = '8.8 °C'
$ws.Range("E25").Value = '2026-02-19 23:19:19'
$ws.Range("I25").Value = '8.0 mm'
$ws.Range("E26").Value = '2026-02-19 23:19:21'
$c = $ws.Range("H26")
$c.NumberFormat = "@"
$c.Value = '55%'
$fmtSrc.Copy()
$c.PasteSpecial(-4122)

$ws.Range("J26").Value = '1010.4 hPa'
$ws.Range("E27").Value = '2026-02-19 23:19:24'
$ws.Range("O27").Value = '-3.6 °C'
$ws.Range("E28").Value = '2026-02-19 23:19:26'
$c = $ws.Range("H28")
$c.NumberFormat = "@"
$c.Value = '66%'
$fmtSrc.Copy()
$c.PasteSpecial(-4122)

$ws.Range("J28").Value = '1010.5 hPa'
$ws.Range("O28").Value = '8.9 °C'
$ws.Range("E29").Value = '2026-02-19 23:19:29'
$ws.Range("N29").Value = '4.0 °C 22:44 TU'
$ws.Range("O29").Value = '9.9 °C'
$ws.Range("E30").Value = '2026-02-19 23:19:31'
$ws.Range("J30").Value = '1010.6 hPa'
$ws.Range("E31").Value = '2026-02-19 23:19:34'
$c = $ws.Range("H31")
$c.NumberFormat = "@"
$c.Value = '49%'
$fmtSrc.Copy()
$c.PasteSpecial(-4122)

$ws.Range("J31").Value = '1009.9 hPa'
$ws.Range("L31").Value = '149.0 km/h - 338º 22:46 TU'
$ws.Range("E32").Value = '2026-02-19 23:19:36'
$ws.Range("O32").Value = '4.6 °C'
$ws.Range("E33").Value = '2026-02-19 23:19:39'
$c = $ws.Range("H33")
$c.NumberFormat = "@"
$c.Value = '58%'
$fmtSrc.Copy()
$c.PasteSpecial(-4122)

$ws.Range("J33").Value = '1011.5 hPa'
$ws.Range("E34").Value = '2026-02-19 23:19:41'
$ws.Range("E35").Value = '2026-02-19 23:19:43'
$c = $ws.Range("H35")
$c.NumberFormat = "@"
$c.Value = '69%'
$fmtSrc.Copy()
$c.PasteSpecial(-4122)

$ws.Range("J35").Value = '1017.1 hPa'
$ws.Range("O35").Value = '3.9 °C'
$ws.Range("E36").Value = '2026-02-19 23:19:46'
$ws.Range("J36").Value = '1010.8 hPa'
$ws.Range("E37").Value = '2026-02-19 23:19:49'
$ws.Range("J37").Value = '1012.0 hPa'
$ws.Range("O37").Value = '5.6 °C'
$ws.Range("E38").Value = '2026-02-19 23:19:51'
$ws.Range("E39").Value = '2026-02-19 23:19:53'
$ws.Range("I39").Value = '5.2 mm'
$ws.Range("E40").Value = '2026-02-19 23:19:56'
$ws.Range("J40").Value = '1013.2 hPa'
$ws.Range("E41").Value = '2026-02-19 23:19:58'
$ws.Range("J41").Value = '1013.5 hPa'
$ws.Range("E42").Value = '2026-02-19 23:20:01'
$c = $ws.Range("H42")
$c.NumberFormat = "@"
$c.Value = '78%'
$fmtSrc.Copy()
$c.PasteSpecial(-4122)

$ws.Range("O42").Value = '10.7 °C'
$ws.Range("E43").Value = '2026-02-19 23:20:03'
$ws.Range("N43").Value = '5.0 °C 22:59 TU'
$ws.Range("E44").Value = '2026-02-19 23:20:05'
$ws.Range("I44").Value = '10.8 mm'
$ws.Range("E45").Value = '2026-02-19 23:20:08'
$c = $ws.Range("H45")
$c.NumberFormat = "@"
$c.Value = '84%'
$fmtSrc.Copy()
$c.PasteSpecial(-4122)

$ws.Range("I45").Value = '3.7 mm'
$ws.Range("J45").Value = '1016.4 hPa'
$ws.Range("E46").Value = '2026-02-19 23:20:10'
$ws.Range("J46").Value = '1016.4 hPa'
$ws.Range("O46").Value = '12.7 °C'

$excel.CutCopyMode = $false
